# DPUB-ARIA suggested-output sheet: post June-1-2021-call terminology tweaks.
#
# The "Suggested Output" column (B) text is revised for five roles
# (doc-backlink, doc-biblioref, doc-glossref, doc-noteref, doc-pullquote)
# and three of those cells also pick up a vertical-center alignment.
# Editing the cells leaves Excel's own "last place you edited" bookmark
# (_GoBack) pointing at the final cell touched (B34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlVAlignCenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# doc-backlink (row 6): "back to referencing item" -> "backlink"
$ws.Range("B6").Value = "backlink"

# doc-biblioref (row 8): "reference to bibliographic entry" -> "Bibliographic reference"
$ws.Range("B8").Value = "Bibliographic reference"
$ws.Range("B8").VerticalAlignment = $xlVAlignCenter

# doc-glossref (row 24): "reference to glossary term" -> "Glossary reference"
$ws.Range("B24").Value = "Glossary reference"
$ws.Range("B24").VerticalAlignment = $xlVAlignCenter

# doc-noteref (row 27): "reference to note item" -> "note reference"
$ws.Range("B27").Value = "note reference"

# doc-pullquote (row 34): "emphasized quote" -> "Emphasised Excerpt"
$ws.Range("B34").Value = "Emphasised Excerpt"
$ws.Range("B34").VerticalAlignment = $xlVAlignCenter

# Excel drops a hidden "_GoBack" bookmark at the last edited cell.
$ws.Names.Add("_GoBack", $ws.Range("B34"))

# Reflect where the author ended up scrolled/selected when they saved.
$ws.Range("B39").Select()
